# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.309.39"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").Value = "3.511.56"
$ws.Range("E3").Value = "  +0.40%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "599.27"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.21%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "174.21"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +3.06%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -1.33%  "
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").Value = "4.114.40"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("E13").Value = "  +0.05%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "30.22"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +7.40%  "
$ws.Range("D15").Value = "67.255.78"
$ws.Range("E15").Value = "  +0.89%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.0000179"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "3.513.21"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("E18").Value = "  -0.21%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "14.55"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +3.55%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "394.51"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.22%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "7.98"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.49%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "73.47"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("E23").Value = "  +0.17%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "0.537"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("E26").Value = "  +0.59%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "10.18"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -0.03%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.997"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  -2.92%  "
$ws.Range("E31").Value = "  -1.78%  "
$ws.Range("E32").Value = "  +0.43%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "23.69"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.43%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "7.39"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +0.72%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.63"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +1.81%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "163.41"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +0.72%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.878"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -2.19%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "1.92"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +0.82%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "6.88"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("E40").Value = "  +0.23%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "27.27"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +1.79%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.0733"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "26.18"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.810.00"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  -0.75%  "
$ws.Range("E47").Value = "  -2.43%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "342.46"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("E49").Value = "  -0.54%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "33.99"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("E51").Value = "  -0.70%  "
